# Generate Report for Handback
# Refresh the report's timestamp fields: the overview sheet's
# "Latest HO Xliff Generate Date" plus each locale sheet's
# "Correspond Handoff Datetime" / "Correspond Handback DateTime".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (row 2, column G)
$wsOverview.Range("G2").Value = "2016-11-03 19:52:42"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-11-03 19:52:28"
$wsZhCn.Range("K2").Value = "2016-11-03 19:53:16"

# de-de sheet: Correspond Handback DateTime (K2)
$wsDeDe.Range("K2").Value = "2016-11-03 19:53:32"
